$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="43.259.33"; E="  -0.90%  "},
    @{B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="2.285.45"; E="  -0.03%  "},
    @{B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.00"; E="  -0.61%  "},
    @{B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="113.58"; E="  +1.07%  "},
    @{B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="266.47"; E="  -0.37%  "},
    @{B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.622"; E="  +0.06%  "},
    @{B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.00"; E="  -0.13%  "},
    @{B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.610"; E="  -0.83%  "},
    @{B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="47.82"; E="  +1.14%  "},
    @{B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.0933"; E="  -0.15%  "},
    @{B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="9.26"; E="  +8.33%  "},
    @{B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.108"; E="  +1.24%  "},
    @{B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="15.58"; E="  +0.20%  "},
    @{B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.628.16"; E="  -0.02%  "},
    @{B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.871"; E="  +2.69%  "},
    @{B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.289.19"; E="  -0.77%  "},
    @{B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="43.275.02"; E="  -0.60%  "},
    @{B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.0000108"; E="  -0.39%  "},
    @{B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="6.85"; E="  +5.35%  "},
    @{B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="71.63"; E="  -0.81%  "},
    @{B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="2.52"; E="  +0.09%  "},
    @{B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="233.03"; E="  +0.21%  "},
    @{B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="9.70"; E="  +2.39%  "},
    @{B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="2.90"; E="  +2.79%  "},
    @{B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.01"; E="  +1.41%  "},
    @{B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="11.41"; E="  +0.27%  "},
    @{B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.99"; E="  +0.26%  "},
    @{B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="40.84"; E="  -5.35%  "},
    @{B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="3.33"; E="  -2.59%  "},
    @{B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.25"; E="  -0.78%  "},
    @{B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="173.45"; E="  -1.47%  "},
    @{B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="21.46"; E="  -0.87%  "},
    @{B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.0909"; E="  -1.50%  "},
    @{B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.78"; E="  +5.62%  "},
    @{B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.127"; E="  +0.70%  "},
    @{B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="4.69"; E="  -0.45%  "},
    @{B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="3.99"; E="  +3.31%  "},
    @{B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0362"; E="  +2.38%  "},
    @{B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.104"; E="  -4.34%  "},
    @{B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.68"; E="  +10.99%  "},
    @{B="MultiversX"; C="https://coinranking.com/coin/omwkOTglq+multiversx-egld"; D="77.89"; E="  +4.11%  "},
    @{B="Celestia"; C="https://coinranking.com/coin/YQcD0lBl7+celestia-tia"; D="14.00"; E="  +5.90%  "},
    @{B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.241"; E="  -0.39%  "},
    @{B="THORChain"; C="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D="6.26"; E="  +5.34%  "},
    @{B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="1.00"; E="  -0.34%  "},
    @{B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.39"; E="  -1.78%  "},
    @{B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="8.72"; E="  -0.23%  "},
    @{B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="104.17"; E="  +2.77%  "},
    @{B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.25"; E="  +1.96%  "},
    @{B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.0996"; E="  -0.31%  "}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $cellD = $ws.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row.D
    $cellD.Style = "Normal"
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}
